$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(206, 1).Value2 = 204
$ws.Cells.Item(206, 2).Value2 = 0.6805098684210524
$ws.Cells.Item(207, 1).Value2 = 205
$ws.Cells.Item(207, 2).Value2 = 0.7033492822966507
$ws.Cells.Item(208, 1).Value2 = 206
$ws.Cells.Item(208, 2).Value2 = 0.75
$ws.Cells.Item(209, 1).Value2 = 207
$ws.Cells.Item(209, 2).Value2 = 0.2690058479532163
$ws.Cells.Item(210, 1).Value2 = 208
$ws.Cells.Item(210, 2).Value2 = 0.5032894736842105
$ws.Cells.Item(211, 1).Value2 = 209
$ws.Cells.Item(211, 2).Value2 = 0.4041353383458646
$ws.Cells.Item(212, 1).Value2 = 210
$ws.Cells.Item(212, 2).Value2 = 0.4725877192982456
$ws.Cells.Item(213, 1).Value2 = 211
$ws.Cells.Item(213, 2).Value2 = 0.4263157894736842
$ws.Cells.Item(214, 1).Value2 = 212
$ws.Cells.Item(214, 2).Value2 = 0.2923519736842105
$ws.Cells.Item(215, 1).Value2 = 213
$ws.Cells.Item(215, 2).Value2 = 0.2321820175438596
$ws.Cells.Item(216, 1).Value2 = 214
$ws.Cells.Item(216, 2).Value2 = 0.4736842105263158
$ws.Cells.Item(217, 1).Value2 = 215
$ws.Cells.Item(217, 2).Value2 = 0.4736842105263158

# Copy formatting (bold, border, center/top alignment) from A205 style to new A cells
$ws.Range("A205").Copy() | Out-Null
$ws.Range("A206:A217").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Write-Output "done"
